# Change the three data tables (slides 14, 15 and 16) from the
# "Table_0" custom table style to the built-in table style
# {520DA83D-3E77-42E7-826F-08A9EFBAAB20}, matching a Table Styles
# gallery change made to those tables in PowerPoint.

$p = $ppt.ActivePresentation

$newStyleId = "{520DA83D-3E77-42E7-826F-08A9EFBAAB20}"
$slideIndexesWithTables = @(14, 15, 16)

foreach ($idx in $slideIndexesWithTables) {
    $slide = $p.Slides.Item($idx)

    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newStyleId)
        }
    }
}
